$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows starting at row 15 (pushes old rows 15-20 down to 20-25,
# so the old "Misc" section that began at row 17 now begins at row 22).
$ws.Rows("15:19").Insert()

# --- New "Power Distribution wire" section -------------------------------

# Section header
$ws.Cells.Item(15, 1).Value2 = "Power Distribution wire"

# Column A (colors) for the three wire rows
$ws.Cells.Item(16, 1).Value2 = "Black"
$ws.Cells.Item(17, 1).Value2 = "White"
$ws.Cells.Item(18, 1).Value2 = "Brown"

# Row 16 (Black wire) - columns B, C, E
$ws.Cells.Item(16, 2).Value2 = "22759/32-12-0"
$ws.Cells.Item(16, 3).Value2 = "22759/32-12-0-DS-ND"
$ws.Cells.Item(16, 5).Value2 = "12 AWG Hook-Up Wire 37/28 Black 600V Enter Number of Feet in Order Quantity"

# Row 17 (White wire) - columns B, C, E
$ws.Cells.Item(17, 2).Value2 = "81044/12-12-9"
$ws.Cells.Item(17, 3).Value2 = "A132407-DS-ND"
$ws.Cells.Item(17, 5).Value2 = "12 AWG Hook-Up, Dual Wall Wire 37/28 White 600V Enter Number of Feet in Order Quantity"

# Row 18 (Brown wire) - columns E, C, B (matches original authoring order)
$ws.Cells.Item(18, 5).Value2 = "12 AWG Hook-Up Wire 37/28 Brown 600V Enter Number of Feet in Order Quantity"
$ws.Cells.Item(18, 3).Value2 = "55A0111-12-1-DS-ND"
$ws.Cells.Item(18, 2).Value2 = "55A0111-12-1"

# --- Fix up the hyperlink that moved from B19 to B24 ----------------------
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B24"), "https://www.digikey.ca/en/products/detail/gc-electronics/CT-4071-1R/6616089", "", "", "https://www.digikey.ca/en/products/detail/gc-electronics/CT-4071-1R/6616089") | Out-Null

# --- Update the active selection ------------------------------------------
$ws.Range("C19:C20").Select()
